# "Removed Test Case Inter-Dependency"
#
# The loan product's name/shortname were hard-coded to values that a
# previous automated test run could already have created in the target
# system, making this test case depend on test execution order/history.
# Here we make the product name/short name unique (independent of any
# prior run) and swap which sheet/cell is active in the workbook.

$wb  = $excel.ActiveWorkbook
$ProductLoanInput  = $wb.Worksheets.Item("ProductLoanInput")
$ProductLoanOutput = $wb.Worksheets.Item("ProductLoanOutput")

# --- content changes -------------------------------------------------

# productname (was "4137-RBI-EPP-FL-SAR-NOREC-MOREREPAY")
$ProductLoanInput.Range("B1").Value  = "4137-RBI-EPP-FL-SAR-NOREC-MOREREPAY-1st"
$ProductLoanOutput.Range("B1").Value = "4137-RBI-EPP-FL-SAR-NOREC-MOREREPAY-1st"

# shortname (was the numeric 4137)
$ProductLoanInput.Range("B2").Value = "413c"

# --- active sheet / selection -----------------------------------------
# Previously "ProductLoanInput" was the selected tab (cell B21 selected)
# and "ProductLoanOutput" was not selected (cell F10 selected). Flip
# that: "ProductLoanOutput" becomes the active tab, and both sheets'
# selections reset to B1.

$ProductLoanInput.Range("B1").Select()
$ProductLoanInput.Activate()

$ProductLoanOutput.Range("B1").Select()
$ProductLoanOutput.Activate()
